$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table by one more year column (O), mirroring the formatting of
# the existing neighbouring cells, then fill in the new 2021 figures.

# Row 2: bottom border strip, same as N2 (blank / border-only cell)
$ws.Range("N2").Copy()
$ws.Range("O2").PasteSpecial(-4122)

# Row 3: year header
$ws.Range("N3").Copy()
$ws.Range("O3").PasteSpecial(-4122)
$ws.Range("O3").Value = 2021

# Row 4: per-person indicator (formula-driven)
$ws.Range("N4").Copy()
$ws.Range("O4").PasteSpecial(-4122)
$ws.Range("O4").Formula = "=O5/O6*1000"

# Row 5: thousand-ton figure (format matches the row's text cells, not N5)
$ws.Range("C5").Copy()
$ws.Range("O5").PasteSpecial(-4122)
$ws.Range("O5").Value = 1229.5999999999999

# Row 6: population figure
$ws.Range("N6").Copy()
$ws.Range("O6").PasteSpecial(-4122)
$ws.Range("O6").Value = 6436.9

$excel.CutCopyMode = $false

# Match the selection left behind in the authored workbook
$ws.Range("P16").Select()
